# Auto-generated Excel COM-interop script
# Updates cached market-price / profit values on the Aegis_Profits (Leve profit) sheets
# to reflect a scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1102.6
$ws.Cells.Item(17, 10).Value = 1102.6
$ws.Cells.Item(17, 12).Value = 3307.8
$ws.Cells.Item(17, 14).Value = -3643.8
$ws.Cells.Item(98, 8).Value = 674.25
$ws.Cells.Item(98, 9).Value = 732.3333
$ws.Cells.Item(98, 10).Value = 500
$ws.Cells.Item(98, 11).Value = 732.3333
$ws.Cells.Item(98, 12).Value = 500
$ws.Cells.Item(98, 13).Value = 765.6667
$ws.Cells.Item(98, 14).Value = -3496
$ws.Cells.Item(103, 8).Value = 600
$ws.Cells.Item(103, 9).Value = 600
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 11).Value = 1800
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 13).Value = -1214
$ws.Cells.Item(122, 8).Value = 674.25
$ws.Cells.Item(122, 9).Value = 732.3333
$ws.Cells.Item(122, 10).Value = 500
$ws.Cells.Item(122, 11).Value = 2196.9999
$ws.Cells.Item(122, 12).Value = 1500
$ws.Cells.Item(122, 13).Value = 253.0001000000002
$ws.Cells.Item(122, 14).Value = -6400
$ws.Cells.Item(131, 8).Value = 4073.68
$ws.Cells.Item(131, 9).Value = 762.3333
$ws.Cells.Item(131, 10).Value = 4525.227
$ws.Cells.Item(131, 11).Value = 2286.9999
$ws.Cells.Item(131, 12).Value = 13575.681
$ws.Cells.Item(131, 13).Value = 2753.0001
$ws.Cells.Item(131, 14).Value = -23655.681
$ws.Cells.Item(137, 8).Value = 1534.2285
$ws.Cells.Item(137, 9).Value = 1406
$ws.Cells.Item(137, 10).Value = 1780
$ws.Cells.Item(137, 11).Value = 4218
$ws.Cells.Item(137, 12).Value = 5340
$ws.Cells.Item(137, 13).Value = -1668
$ws.Cells.Item(137, 14).Value = -10440
$ws.Cells.Item(138, 8).Value = 1883.3334
$ws.Cells.Item(138, 9).Value = 1463.1305
$ws.Cells.Item(138, 11).Value = 4389.3915
$ws.Cells.Item(138, 13).Value = 750.6085000000003
$ws.Range("N103").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 185.33333
$ws.Cells.Item(4, 9).Value = 100
$ws.Cells.Item(4, 11).Value = 100
$ws.Cells.Item(4, 13).Value = 16
$ws.Cells.Item(32, 8).Value = 27866.084
$ws.Cells.Item(32, 9).Value = 5741.875
$ws.Cells.Item(32, 10).Value = 230144.58
$ws.Cells.Item(32, 11).Value = 5741.875
$ws.Cells.Item(32, 12).Value = 230144.58
$ws.Cells.Item(32, 13).Value = -5454.875
$ws.Cells.Item(32, 14).Value = -230718.58
$ws.Cells.Item(61, 8).Value = 2198.625
$ws.Cells.Item(61, 9).Value = 2198.625
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 2198.625
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -1986.625
$ws.Cells.Item(74, 8).Value = 1324.24
$ws.Cells.Item(74, 9).Value = 1439.1765
$ws.Cells.Item(74, 10).Value = 1080
$ws.Cells.Item(74, 11).Value = 1439.1765
$ws.Cells.Item(74, 12).Value = 1080
$ws.Cells.Item(74, 13).Value = -565.1765
$ws.Cells.Item(74, 14).Value = -2828
$ws.Cells.Item(77, 8).Value = 1324.24
$ws.Cells.Item(77, 9).Value = 1439.1765
$ws.Cells.Item(77, 10).Value = 1080
$ws.Cells.Item(77, 11).Value = 7195.8825
$ws.Cells.Item(77, 12).Value = 5400
$ws.Cells.Item(77, 13).Value = -2827.8825
$ws.Cells.Item(77, 14).Value = -14136
$ws.Cells.Item(102, 8).Value = 45145.566
$ws.Cells.Item(102, 9).Value = 101647
$ws.Cells.Item(102, 10).Value = 1682.9231
$ws.Cells.Item(102, 11).Value = 101647
$ws.Cells.Item(102, 12).Value = 1682.9231
$ws.Cells.Item(102, 13).Value = -100025
$ws.Cells.Item(102, 14).Value = -4926.9231
$ws.Cells.Item(122, 8).Value = 2367.611
$ws.Cells.Item(122, 9).Value = 1630.7693
$ws.Cells.Item(122, 11).Value = 4892.3079
$ws.Cells.Item(122, 13).Value = -2442.3079
$ws.Cells.Item(132, 8).Value = 1945.4894
$ws.Cells.Item(132, 9).Value = 1522.7435
$ws.Cells.Item(132, 10).Value = 4006.375
$ws.Cells.Item(132, 11).Value = 4568.2305
$ws.Cells.Item(132, 12).Value = 12019.125
$ws.Cells.Item(132, 13).Value = -2038.2305
$ws.Cells.Item(132, 14).Value = -17079.125
$ws.Cells.Item(136, 8).Value = 2198.625
$ws.Cells.Item(136, 9).Value = 2198.625
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 6595.875
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -4045.875
$ws.Range("N136").ClearContents()
$ws.Range("N61").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 34086.516
$ws.Cells.Item(105, 9).Value = 38817.965
$ws.Cells.Item(105, 10).Value = 2149.25
$ws.Cells.Item(105, 11).Value = 38817.965
$ws.Cells.Item(105, 12).Value = 2149.25
$ws.Cells.Item(105, 13).Value = -37070.965
$ws.Cells.Item(105, 14).Value = -5643.25
$ws.Cells.Item(133, 8).Value = 359999
$ws.Cells.Item(133, 10).Value = 359999
$ws.Cells.Item(133, 12).Value = 359999
$ws.Cells.Item(134, 8).Value = 1907.0294
$ws.Cells.Item(134, 9).Value = 1825.3438
$ws.Cells.Item(134, 10).Value = 3214
$ws.Cells.Item(134, 11).Value = 5476.0314
$ws.Cells.Item(134, 12).Value = 9642
$ws.Cells.Item(134, 13).Value = -2941.0314
$ws.Cells.Item(133, 14).Value = -370119
$ws.Cells.Item(134, 14).Value = -14712

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1454.0714
$ws.Cells.Item(58, 9).Value = 1221.7391
$ws.Cells.Item(58, 10).Value = 2522.8
$ws.Cells.Item(58, 11).Value = 1221.7391
$ws.Cells.Item(58, 12).Value = 2522.8
$ws.Cells.Item(58, 13).Value = -1018.7391
$ws.Cells.Item(58, 14).Value = -2928.8
$ws.Cells.Item(105, 8).Value = 1162.2941
$ws.Cells.Item(105, 9).Value = 1056
$ws.Cells.Item(105, 10).Value = 1658.3334
$ws.Cells.Item(105, 11).Value = 1056
$ws.Cells.Item(105, 12).Value = 1658.3334
$ws.Cells.Item(105, 13).Value = 691
$ws.Cells.Item(105, 14).Value = -5152.3334
$ws.Cells.Item(122, 8).Value = 553.125
$ws.Cells.Item(122, 9).Value = 348.33334
$ws.Cells.Item(122, 10).Value = 676
$ws.Cells.Item(122, 11).Value = 1045.00002
$ws.Cells.Item(122, 12).Value = 2028
$ws.Cells.Item(122, 13).Value = 1404.99998
$ws.Cells.Item(122, 14).Value = -6928
$ws.Cells.Item(132, 8).Value = 5192.5264
$ws.Cells.Item(132, 9).Value = 4918.5713
$ws.Cells.Item(132, 11).Value = 14755.7139
$ws.Cells.Item(132, 13).Value = -12225.7139
$ws.Cells.Item(134, 8).Value = 804.63635
$ws.Cells.Item(134, 9).Value = 840.1
$ws.Cells.Item(134, 10).Value = 450
$ws.Cells.Item(134, 11).Value = 2520.3
$ws.Cells.Item(134, 12).Value = 1350
$ws.Cells.Item(134, 13).Value = 14.69999999999982
$ws.Cells.Item(134, 14).Value = -6420
$ws.Cells.Item(136, 8).Value = 1454.0714
$ws.Cells.Item(136, 9).Value = 1221.7391
$ws.Cells.Item(136, 10).Value = 2522.8
$ws.Cells.Item(136, 11).Value = 3665.2173
$ws.Cells.Item(136, 12).Value = 7568.400000000001
$ws.Cells.Item(136, 13).Value = -1115.2173
$ws.Cells.Item(136, 14).Value = -12668.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 1221.7778
$ws.Cells.Item(75, 10).Value = 1349.3334
$ws.Cells.Item(75, 12).Value = 4048.0002
$ws.Cells.Item(75, 14).Value = -6044.0002
$ws.Cells.Item(78, 8).Value = 1221.7778
$ws.Cells.Item(78, 10).Value = 1349.3334
$ws.Cells.Item(78, 12).Value = 12144.0006
$ws.Cells.Item(78, 14).Value = -22128.0006
$ws.Cells.Item(81, 8).Value = 3395.5715
$ws.Cells.Item(81, 9).Value = 887
$ws.Cells.Item(81, 10).Value = 4399
$ws.Cells.Item(81, 11).Value = 2661
$ws.Cells.Item(81, 12).Value = 13197
$ws.Cells.Item(81, 13).Value = -1538
$ws.Cells.Item(81, 14).Value = -15443
$ws.Cells.Item(84, 8).Value = 3395.5715
$ws.Cells.Item(84, 9).Value = 887
$ws.Cells.Item(84, 10).Value = 4399
$ws.Cells.Item(84, 11).Value = 7983
$ws.Cells.Item(84, 12).Value = 39591
$ws.Cells.Item(84, 13).Value = -2367
$ws.Cells.Item(84, 14).Value = -50823
$ws.Cells.Item(107, 8).Value = 463862.97
$ws.Cells.Item(107, 9).Value = 664.44446
$ws.Cells.Item(107, 10).Value = 662376.6
$ws.Cells.Item(107, 11).Value = 1993.33338
$ws.Cells.Item(107, 12).Value = 1987129.8
$ws.Cells.Item(107, 13).Value = -73.33338000000003
$ws.Cells.Item(107, 14).Value = -1990969.8
$ws.Cells.Item(131, 8).Value = 6708.0957
$ws.Cells.Item(131, 10).Value = 6708.0957
$ws.Cells.Item(131, 12).Value = 20124.2871
$ws.Cells.Item(131, 14).Value = -30204.2871

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3648.6667
$ws.Cells.Item(126, 10).Value = 1975
$ws.Cells.Item(126, 12).Value = 5925
$ws.Cells.Item(126, 14).Value = -10865
$ws.Cells.Item(132, 8).Value = 2634.0356
$ws.Cells.Item(132, 9).Value = 2465.5
$ws.Cells.Item(132, 10).Value = 3252
$ws.Cells.Item(132, 11).Value = 7396.5
$ws.Cells.Item(132, 12).Value = 9756
$ws.Cells.Item(132, 13).Value = -4866.5
$ws.Cells.Item(132, 14).Value = -14816

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4323.75
$ws.Cells.Item(7, 9).Value = 3000
$ws.Cells.Item(7, 10).Value = 4512.857
$ws.Cells.Item(7, 11).Value = 3000
$ws.Cells.Item(7, 12).Value = 4512.857
$ws.Cells.Item(7, 13).Value = -2888
$ws.Cells.Item(7, 14).Value = -4736.857
$ws.Cells.Item(40, 8).Value = 51985.45
$ws.Cells.Item(40, 9).Value = 112501.336
$ws.Cells.Item(40, 11).Value = 112501.336
$ws.Cells.Item(40, 13).Value = -112365.336
$ws.Cells.Item(55, 8).Value = 392758.03
$ws.Cells.Item(55, 9).Value = 758700.9399999999
$ws.Cells.Item(55, 10).Value = 676.3570999999999
$ws.Cells.Item(55, 11).Value = 758700.9399999999
$ws.Cells.Item(55, 12).Value = 676.3570999999999
$ws.Cells.Item(55, 13).Value = -758527.9399999999
$ws.Cells.Item(55, 14).Value = -1022.3571
$ws.Cells.Item(126, 8).Value = 4323.75
$ws.Cells.Item(126, 9).Value = 3000
$ws.Cells.Item(126, 10).Value = 4512.857
$ws.Cells.Item(126, 11).Value = 9000
$ws.Cells.Item(126, 12).Value = 13538.571
$ws.Cells.Item(126, 13).Value = -6530
$ws.Cells.Item(126, 14).Value = -18478.571
$ws.Cells.Item(132, 8).Value = 3382
$ws.Cells.Item(132, 9).Value = 3271.8438
$ws.Cells.Item(132, 10).Value = 3885.5715
$ws.Cells.Item(132, 11).Value = 9815.5314
$ws.Cells.Item(132, 12).Value = 11656.7145
$ws.Cells.Item(132, 13).Value = -7285.5314
$ws.Cells.Item(132, 14).Value = -16716.7145
$ws.Cells.Item(136, 8).Value = 1405.4
$ws.Cells.Item(136, 9).Value = 1269.2812
$ws.Cells.Item(136, 10).Value = 1949.875
$ws.Cells.Item(136, 11).Value = 3807.8436
$ws.Cells.Item(136, 12).Value = 5849.625
$ws.Cells.Item(136, 13).Value = -1257.8436
$ws.Cells.Item(136, 14).Value = -10949.625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2383.9375
$ws.Cells.Item(122, 9).Value = 1699.6
$ws.Cells.Item(122, 10).Value = 2695
$ws.Cells.Item(122, 11).Value = 5098.799999999999
$ws.Cells.Item(122, 12).Value = 8085
$ws.Cells.Item(122, 13).Value = -2648.799999999999
$ws.Cells.Item(122, 14).Value = -12985
$ws.Cells.Item(126, 8).Value = 1364
$ws.Cells.Item(126, 9).Value = 1248.7778
$ws.Cells.Item(126, 10).Value = 1623.25
$ws.Cells.Item(126, 11).Value = 3746.3334
$ws.Cells.Item(126, 12).Value = 4869.75
$ws.Cells.Item(126, 13).Value = -1276.3334
$ws.Cells.Item(126, 14).Value = -9809.75
$ws.Cells.Item(132, 8).Value = 2025.2554
$ws.Cells.Item(132, 9).Value = 1946.0731
$ws.Cells.Item(132, 10).Value = 2566.3333
$ws.Cells.Item(132, 11).Value = 5838.219300000001
$ws.Cells.Item(132, 12).Value = 7698.999899999999
$ws.Cells.Item(132, 13).Value = -3308.219300000001
$ws.Cells.Item(132, 14).Value = -12758.9999
$ws.Cells.Item(136, 8).Value = 2672.4375
$ws.Cells.Item(136, 9).Value = 868.5454999999999
$ws.Cells.Item(136, 10).Value = 6641
$ws.Cells.Item(136, 11).Value = 2605.6365
$ws.Cells.Item(136, 12).Value = 19923
$ws.Cells.Item(136, 13).Value = -55.63649999999961
$ws.Cells.Item(136, 14).Value = -25023
